$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 1.03
$ws.Range("U2").Value = 1.03

# Row 4
$ws.Range("F4").Value = 1.58
$ws.Range("G4").Value = 1.65
$ws.Range("H4").Value = 6.2
$ws.Range("I4").Value = 7.2
$ws.Range("J4").Value = 4.2
$ws.Range("K4").Value = 4.7
$ws.Range("P4").Value = 2

# Row 6
$ws.Range("F6").Value = 1.67
$ws.Range("G6").Value = 980
$ws.Range("H6").Value = 1.04
$ws.Range("I6").Value = 9.4
$ws.Range("J6").Value = 3.4

# Row 7
$ws.Range("J7").Value = 3.05

# Row 8
$ws.Range("G8").Value = 3.4
$ws.Range("H8").Value = 2.12
$ws.Range("I8").Value = 3.35

# Row 9
$ws.Range("H9").Value = 4.3

# Row 10
$ws.Range("L10").Value = 1.48
$ws.Range("AN10").Value = 55

# Row 11
$ws.Range("F11").Value = 2.6
$ws.Range("G11").Value = 3.3
$ws.Range("H11").Value = 2.8
$ws.Range("I11").Value = 3.6
$ws.Range("J11").Value = 2.72
$ws.Range("K11").Value = 3.5
$ws.Range("P11").Value = 1.58

# Row 12
$ws.Range("J12").Value = 2.64

# Row 13
$ws.Range("D13").Value = "CSD Rangers"
$ws.Range("E13").Value = "Antofagasta"
$ws.Range("F13").Value = 1.04
$ws.Range("G13").Value = 1000
$ws.Range("H13").Value = 1.04
$ws.Range("J13").Value = 2.92
$ws.Range("P13").Value = 1.63
$ws.Range("Q13").Value = 1.9

# Row 14
$ws.Range("D14").Value = "San Luis"
$ws.Range("E14").Value = "Deportes Recoleta"
$ws.Range("F14").Value = 2.18
$ws.Range("G14").Value = 3.05
$ws.Range("H14").Value = 1.49
$ws.Range("I14").Value = 5.7
$ws.Range("J14").Value = 2.74
$ws.Range("P14").Value = 1.53
$ws.Range("Q14").Value = 2.04
